# Apply the row-10/11/13 data rotation described by the commit diff.
# (Row 10's old data -> row 11, row 11's old data -> row 13, row 13's old data -> row 10,
#  plus the two "Observatörer" (AX) tweaks that came along with the re-ordering.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $value) {
    # Force the cell to keep a TEXT value even when it looks like a number
    # or a date (e.g. "1", "3", "2023-08-12"), then restore the default
    # "Normal" style so no stray number-format sticks around.
    $cell = $ws.Range($rng)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---- Row 10 (becomes the old row 13 data) ----
$ws.Range("A10").Value = 111611158
$ws.Range("B10").Value = 86021
$ws.Range("E10").Value = 4037
$ws.Range("F10").Value = "Bolmörtsskivling"
$ws.Range("G10").Value = "Entoloma sinuatum"
$ws.Range("H10").Value = "(Bull.) P.Kumm."
Set-TextValue "I10" "1"
$ws.Range("Q10").Value = 663128.0992466732
$ws.Range("R10").Value = 6634761.25188593
$ws.Range("AC10").Value = "1 ex. i lövförna under ek och hassel."

# ---- Row 11 (becomes the old row 10 data) ----
$ws.Range("A11").Value = 111611138
$ws.Range("B11").Value = 81796
$ws.Range("E11").Value = 5406
$ws.Range("F11").Value = "Gulmjölkig storskål"
$ws.Range("G11").Value = "Peziza succosa"
$ws.Range("H11").Value = "Berk."
Set-TextValue "I11" "3"
$ws.Range("Q11").Value = 663213.3366271106
$ws.Range("R11").Value = 6634830.464506784
Set-TextValue "Y11" "2023-08-12"
Set-TextValue "AA11" "2023-08-12"
$ws.Range("AC11").Value = "3 ex. på bar jord och i lövförna."
$ws.Range("AX11").Value = "Gillis Aronsson, Cajsa Björkén"

# ---- Row 13 (becomes the old row 11 data) ----
$ws.Range("A13").Value = 111611146
$ws.Range("B13").Value = 88630
$ws.Range("E13").Value = 4823
$ws.Range("F13").Value = "Hasselsopp"
$ws.Range("G13").Value = "Leccinellum pseudoscabrum"
$ws.Range("H13").Value = "(Kallenb.) Mikšík"
$ws.Range("Q13").Value = 663088.0668624006
$ws.Range("R13").Value = 6634684.960451891
Set-TextValue "Y13" "2023-08-11"
Set-TextValue "AA13" "2023-08-11"
$ws.Range("AC13").Value = "1 ex. under ek och hassel."
$ws.Range("AX13").Value = "Gillis Aronsson"
